$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph (2nd paragraph of the doc)
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
$metaRange.Delete()

# ---------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Drago: Jewels of Fortune for
#    Free - Review" right before the last paragraph (the one that used
#    to read "Please create an image for ...").
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($count - 1)
$secondToLast.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($count)
$newPara.Format.Style = $d.Styles.Item("Normal")

$newParaRange = $newPara.Range
$insertPoint = $d.Range($newParaRange.End - 1, $newParaRange.End - 1)
$insertPoint.InsertAfter("Play Drago: Jewels of Fortune for Free - Review")
$insertPoint.Bold = $true

# ---------------------------------------------------------------------
# 3) Replace the text of the final paragraph (still italic) from the
#    image-generation prompt to the new meta-description sentence.
# ---------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "Please create an image for " + [char]34 + "Drago: Jewels of Fortune" + [char]34 + " featuring a happy Maya warrior with glasses in a cartoon style. The warrior should be depicted in a jungle setting with treasure chests and dragons in the background. The image should be vibrant and colorful, with attention to detail in the warrior's clothing and accessories. The overall vibe should be adventurous and exciting, reflecting the theme of the game. Thank you!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Try out Drago: Jewels of Fortune for free and read a comprehensive review of its betting options, features, and accessibility on different devices.",
    2)
